$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Escenas PencilTest")

# 1. Fix Frames for escena_04 (row 5) - cleanup frames 21-41
$ws.Range("C5").Value = 41

# 2. Add the "CleanUp" column to the table
$lo = $ws.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()
$ws.Range("H1").Value = "CleanUp"

# 3. Format header cell H1 like the other header cells (border + center/middle)
$h1 = $ws.Range("H1")
$b1 = $ws.Range("B1")
$h1.HorizontalAlignment = $b1.HorizontalAlignment
$h1.VerticalAlignment = $b1.VerticalAlignment
$h1.Borders.Item(7).LineStyle = $b1.Borders.Item(7).LineStyle
$h1.Borders.Item(10).LineStyle = $b1.Borders.Item(10).LineStyle
$h1.Borders.Item(8).LineStyle = $b1.Borders.Item(8).LineStyle
$h1.Borders.Item(9).LineStyle = $b1.Borders.Item(9).LineStyle
$h1.Interior.Pattern = -4142

# 4. Format body cells H2:H81 like G2:G81 (same look as rest of table)
$g = $ws.Range("G2:G81")
$g.Copy()
$h = $ws.Range("H2:H81")
$h.PasteSpecial(-4122)

# 5. Mark row 31 (scene 30) as cleaned up
$ws.Range("H31").Value = "X"

# 6. Format totals row cell H82 (bold, centered, no border) like the rest of the totals row
$g82 = $ws.Range("G82")
$g82.Copy()
$h82 = $ws.Range("H82")
$h82.PasteSpecial(-4122)
$h82.Borders.Item(7).LineStyle = -4142
$h82.Borders.Item(8).LineStyle = -4142
$h82.Borders.Item(9).LineStyle = -4142
$h82.Borders.Item(10).LineStyle = -4142

$excel.CutCopyMode = 0

# 7. Conditional formatting: highlight non-blank CleanUp marks
$cf = $ws.Range("H2:H81").FormatConditions.Add(2, 0, "LEN(TRIM(H2))>0")
$cf.Interior.Color = 52377
$cf.Font.Color = 10284544

Write-Output "Done"
